$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")
$ws.Activate()

# Remove the "Number of VL tests recommended per person per year" row (row 41)
$ws.Rows("41").Delete()

$ws.Range("A41").Select()
